$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("metadata")

# Insert a new row above row 15 (shifts CAPTURE DATE related rows, SAMPLING_PLATFORM, etc. down by one)
$ws.Rows.Item(15).Insert()

# Fill the new row with the CAPTURE_QUARTER field definition
$ws.Range("A15").Value = "CAPTURE_QUARTER"
$ws.Range("B15").Value = 'Quarter of the "average" date of capture'

# Update the active selection to match the final saved view
$ws.Range("B6").Select()
